$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''37.822.25'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').Value = '''2.080.33'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''  -0.19%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = '''  +0.02%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = '''233.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '''  -0.26%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('E6').Value = '''  +0.00%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = '''58.55'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '''  -0.62%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('E8').Value = '''  +0.01%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = '''0.393'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''  +0.60%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = '''0.0786'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''  -0.64%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('E11').Value = '''  +3.41%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = '''15.00'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''  +1.54%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = '''2.387.45'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''  -0.20%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = '''21.30'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''  +0.18%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''0.780'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''  +1.36%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''5.39'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''  +1.65%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''2.079.88'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''  -0.22%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = '''37.786.58'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '''  +0.02%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').Value = '''  -1.19%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''71.35'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''  +0.01%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = '''0.0₃0839'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '''  +0.53%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''230.45'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''  +0.64%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').Value = '''  -0.13%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('E24').Value = '''  -0.64%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''2.40'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''  +1.11%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''9.85'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''  +9.22%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''172.04'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '''  +0.94%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('E28').Value = '''  -2.11%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''19.51'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '''  -0.19%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '''1.41'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '''  -0.42%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('E31').Value = '''  +1.37%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = '''4.72'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '''  +0.29%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = '''  +0.62%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = '''4.69'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '''  -0.45%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = '''2.48'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '''  -1.44%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = '''  -0.71%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('E37').Value = '''  -2.21%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('E38').Value = '''  +0.00%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = '''5.48'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '''  +0.75%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = '''0.0236'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''  +9.95%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = '''102.64'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''  +4.01%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('E42').Value = '''  -2.01%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = '''2.93'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''  -0.63%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''16.80'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''  +4.27%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = '''1.451.44'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''  -1.01%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('E46').Value = '''  -1.41%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').Value = '''  -0.69%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = '''  -7.36%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('E50').Value = '''  -1.44%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = '''2.272.40'
$ws.Range('D51').Style = "Normal"
